# Rename the embedded logo pictures that live in the document's header/footer
# stories:
#   - Pearson logo inline pictures (wp:docPr id="2" / id="3", currently
#     name="image2.png") -> name="image1.png"
#   - BTEC logo inline picture   (wp:docPr id="1", currently name="image1.jpg")
#     -> name="image2.jpg"
#
# InlineShape has no writable "Name" in the Word object model (only Shape /
# ShapeRange expose it), so each inline picture is promoted to a floating
# Shape, renamed there, then converted back to an inline shape so the
# <wp:inline> wrapper (and everything else about it) is preserved.

function Rename-InlinePicture($inlineShape, $newName) {
    $floating = $inlineShape.ConvertToShape()
    $floating.Name = $newName
    [void]$floating.ConvertToInlineShape()
}

$pearsonDescr = "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png"
$btecDescr = "BTec_Logo-Orange"

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    # Headers: rename the BTEC logo wherever it shows up.
    for ($hi = 1; $hi -le $sec.Headers.Count; $hi++) {
        $hdr = $sec.Headers.Item($hi)
        if (-not $hdr.Exists) { continue }
        for ($si = 1; $si -le $hdr.Range.InlineShapes.Count; $si++) {
            $shp = $hdr.Range.InlineShapes.Item($si)
            if ($shp.AlternativeText -eq $btecDescr) {
                Rename-InlinePicture $shp "image2.jpg"
            }
        }
    }

    # Footers: rename the Pearson logo wherever it shows up.
    for ($fi = 1; $fi -le $sec.Footers.Count; $fi++) {
        $ftr = $sec.Footers.Item($fi)
        if (-not $ftr.Exists) { continue }
        for ($si = 1; $si -le $ftr.Range.InlineShapes.Count; $si++) {
            $shp = $ftr.Range.InlineShapes.Item($si)
            if ($shp.AlternativeText -eq $pearsonDescr) {
                Rename-InlinePicture $shp "image1.png"
            }
        }
    }
}
